$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 585
$ws.Cells.Item(5, 6).Value = 2570
$ws.Cells.Item(10, 6).Value = 5271
$ws.Cells.Item(11, 6).Value = 98
$ws.Cells.Item(12, 6).Value = 1461
$ws.Cells.Item(14, 6).Value = 591
$ws.Cells.Item(15, 6).Value = 6957
$ws.Cells.Item(16, 6).Value = 388
$ws.Cells.Item(17, 6).Value = 47
$ws.Cells.Item(19, 6).Value = 67
$ws.Cells.Item(20, 6).Value = 4667
$ws.Cells.Item(22, 6).Value = 70
$ws.Cells.Item(23, 6).Value = 2323
$ws.Cells.Item(24, 6).Value = 1252
$ws.Cells.Item(25, 6).Value = 441
$ws.Cells.Item(26, 6).Value = 1153
$ws.Cells.Item(27, 6).Value = 215
$ws.Cells.Item(30, 6).Value = 158
$ws.Cells.Item(31, 6).Value = 365
$ws.Cells.Item(32, 6).Value = 1272
$ws.Cells.Item(33, 6).Value = 1989
$ws.Cells.Item(34, 6).Value = 232
$ws.Cells.Item(35, 6).Value = 515
$ws.Cells.Item(36, 6).Value = 201
$ws.Cells.Item(37, 6).Value = 1371
$ws.Cells.Item(39, 6).Value = 87
$ws.Cells.Item(40, 6).Value = 518
$ws.Cells.Item(41, 6).Value = 164
$ws.Cells.Item(42, 6).Value = 1115
$ws.Cells.Item(43, 6).Value = 2401
$ws.Cells.Item(45, 6).Value = 65
$ws.Cells.Item(47, 6).Value = 231
$ws.Cells.Item(48, 6).Value = 70
$ws.Cells.Item(49, 6).Value = 15

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 426
$ws.Cells.Item(12, 6).Value = 381
$ws.Cells.Item(13, 6).Value = 269
$ws.Cells.Item(15, 6).Value = 39
$ws.Cells.Item(28, 6).Value = 281

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 493
$ws.Cells.Item(6, 6).Value = 1652
$ws.Cells.Item(8, 6).Value = 1282
$ws.Cells.Item(10, 6).Value = 1735
$ws.Cells.Item(11, 6).Value = 2174
$ws.Cells.Item(12, 6).Value = 606
$ws.Cells.Item(13, 6).Value = 510

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 493
$ws.Cells.Item(3, 6).Value = 1652
$ws.Cells.Item(5, 6).Value = 585
$ws.Cells.Item(7, 6).Value = 2570
$ws.Cells.Item(9, 6).Value = 1282
$ws.Cells.Item(10, 6).Value = 2175
$ws.Cells.Item(11, 6).Value = 5271
$ws.Cells.Item(12, 6).Value = 606
$ws.Cells.Item(15, 6).Value = 98
$ws.Cells.Item(17, 6).Value = 1461
$ws.Cells.Item(19, 6).Value = 591
$ws.Cells.Item(20, 6).Value = 6957
$ws.Cells.Item(21, 6).Value = 388
$ws.Cells.Item(22, 6).Value = 510
$ws.Cells.Item(23, 6).Value = 47
$ws.Cells.Item(24, 6).Value = 4667
$ws.Cells.Item(25, 6).Value = 2323
$ws.Cells.Item(26, 6).Value = 1252
$ws.Cells.Item(27, 6).Value = 441
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 215
$ws.Cells.Item(31, 6).Value = 269
$ws.Cells.Item(33, 6).Value = 158
$ws.Cells.Item(35, 6).Value = 365
$ws.Cells.Item(36, 6).Value = 1989
$ws.Cells.Item(37, 6).Value = 232
$ws.Cells.Item(38, 6).Value = 515
$ws.Cells.Item(40, 6).Value = 1371
$ws.Cells.Item(42, 6).Value = 164
$ws.Cells.Item(44, 6).Value = 1115
$ws.Cells.Item(45, 6).Value = 2401
$ws.Cells.Item(46, 6).Value = 65
$ws.Cells.Item(47, 6).Value = 231
$ws.Cells.Item(48, 6).Value = 70
